# QB siteyml change 1/21
# Rename the "Occurrence" column (E) to "Dates Used " and replace the
# placeholder occurrence-number lists in column E with the actual date
# ranges used for each survey wave.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$datesAll = "Aug 24 - Aug Oct 5 - Oct 8, Sept 8 - Sept 10, Sept 21 - Sept 24, Oct 5 - Oct 8, Nov 3 - Nov 5, Nov 17 - Nov 19, 35, Dec 15 - Dec 17, Jan 5 - Jan 7, Feb 3 - Feb 6, March 2 - March 5, March 31 - April 2, April 28 - May 1, May 25 - May 28, June 23 - June 25, July 21 - July 23, Aug 4 - Aug 6, September 15 - Septeber 17, November 10 - November 16, December 8 - December 14, January 11 - January 18, January 3 - January 24"
$datesE4 = "Aug 24 - Aug Oct 5 - Oct 8, Sept 8 - Sept 10, Sept 21 - Sept 24, Oct 5 - Oct 8, Nov 3 - Nov 5, Nov 17 - Nov 19, 35, Dec 15 - Dec 17, Jan 5 - Jan 7, Feb 3 - Feb 6, March 2 - March 5, March 31 - April 2, April 28 - May 1, May 25 - May 28, June 23 - June 25, July 21 - July 23, Aug 4 - Aug 6, September 15 - Septeber 17, November 10 - November 16 , January 11 - January 18, January 3 - January 24"

# Header
$ws.Range("E1").Value = "Dates Used "

# Data rows
$ws.Range("E2").Value = $datesAll
$ws.Range("E3").Value = $datesAll
$ws.Range("E4").Value = $datesE4
$ws.Range("E5").Value = $datesAll
$ws.Range("E6").Value = $datesAll
$ws.Range("E7").Value = $datesAll

# Match the author's final selection: the whole of column E, active cell E1
$ws.Columns("E:E").Select()
